# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (copied from "2022-Q3" as a template,
#    keeping the same fund code/name) right after "总计", and update its
#    quarterly figures.
# 2. Insert a new row into the "总计" summary sheet for the 2022-Q4 period,
#    shifting the existing rows down and renumbering the index column.
# 3. Restore "2020-Q4" (the last sheet) as the active/selected tab, matching
#    the workbook's original view state.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# --- Step 1: create the new "2022-Q4" sheet -------------------------------
$wsQ3.Copy($null, $wsTotal)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

$wsQ4.Range("D2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "0.43"
$wsQ4.Range("E2").NumberFormat = "@"
$wsQ4.Range("E2").Value = "92.90"
$wsQ4.Range("F2").NumberFormat = "@"
$wsQ4.Range("F2").Value = "4.27"
$wsQ4.Range("G2").NumberFormat = "@"
$wsQ4.Range("G2").Value = "0.0184"

# --- Step 2: update the "总计" (totals) summary sheet ----------------------
$wsTotal.Rows.Item(2).Insert()

# Copy the numeric index column's formatting down into the newly opened row
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
# The rest of the new row should look like the plain (unstyled) data rows
$wsTotal.Range("B2:D2").Style = "Normal"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.02

# Renumber the index column for the rows that shifted down
for ($r = 3; $r -le 10; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# --- Step 3: keep "2020-Q4" as the active tab ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$lastSheet.Activate()
